# Applies the cryptos-list price/volume refresh for rows 2-51.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.742.26"
$ws.Range("E2").Value = "  -0.51%  "

$ws.Range("D3").Value = "2.268.35"
$ws.Range("E3").Value = "  -0.64%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'248.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.08%  "

$ws.Range("E6").Value = "  +0.38%  "

$ws.Range("D7").Value = "'76.78"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +6.48%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").Value = "'0.650"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.66%  "

$ws.Range("D10").Value = "'39.70"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.19%  "

$ws.Range("D11").Value = "'0.0968"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.75%  "

$ws.Range("E12").Value = "  -2.49%  "

$ws.Range("E13").Value = "  -0.68%  "

$ws.Range("D14").Value = "2.607.13"
$ws.Range("E14").Value = "  -0.64%  "

$ws.Range("E15").Value = "  -0.12%  "

$ws.Range("D16").Value = "'0.862"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.32%  "

$ws.Range("D17").Value = "2.262.12"
$ws.Range("E17").Value = "  +0.36%  "

$ws.Range("D18").Value = "42.638.91"
$ws.Range("E18").Value = "  -0.60%  "

$ws.Range("D19").Value = "0.0₃0988"
$ws.Range("E19").Value = "  -3.07%  "

$ws.Range("E20").Value = "  -2.37%  "

$ws.Range("D21").Value = "'71.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.31%  "

$ws.Range("D22").Value = "'232.65"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.85%  "

$ws.Range("E23").Value = "  -2.41%  "

$ws.Range("E24").Value = "  -6.10%  "

$ws.Range("E25").Value = "  +0.01%  "

$ws.Range("D26").Value = "'11.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.65%  "

$ws.Range("D27").Value = "'2.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.08%  "

$ws.Range("E28").Value = "  +1.81%  "

$ws.Range("D29").Value = "'167.80"
$ws.Range("D29").Style = "Normal"

$ws.Range("D30").Value = "'20.84"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.30%  "

$ws.Range("E31").Value = "  -2.97%  "

$ws.Range("D32").Value = "'0.0851"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.97%  "

$ws.Range("E33").Value = "  -3.64%  "

$ws.Range("D34").Value = "'30.33"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.44%  "

$ws.Range("E35").Value = "  +0.46%  "

$ws.Range("D36").Value = "'4.54"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.62%  "

$ws.Range("E37").Value = "  -2.16%  "

$ws.Range("E38").Value = "  -3.47%  "

$ws.Range("D39").Value = "'13.84"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.07%  "

$ws.Range("E40").Value = "  -3.47%  "

$ws.Range("D41").Value = "'5.81"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.81%  "

$ws.Range("D42").Value = "'0.208"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.17%  "

$ws.Range("E43").Value = "  +14.83%  "

$ws.Range("D44").Value = "'60.94"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.62%  "

$ws.Range("D45").Value = "'8.84"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.08%  "

$ws.Range("E46").Value = "  -2.21%  "

$ws.Range("E47").Value = "  -0.14%  "

$ws.Range("D48").Value = "'4.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -10.56%  "

$ws.Range("D49").Value = "'1.14"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.69%  "

$ws.Range("E50").Value = "  -3.61%  "

$ws.Range("E51").Value = "  -1.04%  "
